$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that contained "NaN" placeholder values
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()

$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()

# Update selection to match the final state (cell G6 selected)
$ws.Range("G6").Select()
